$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elem = $wb.Worksheets.Item("Elements")

# --- Remove the duplicated "Contact" row (old row 11) so everything below shifts up ---
# (row 10 held "Contact" / "No display for ContactDetail" twice; deleting row 11
#  collapses the table from 21 rows to 20 and shifts rows 12..21 up to 11..20)
$meta.Rows.Item(11).Delete()

# --- Version bump ---
$meta.Range("B3").Value = "6.0.0"

# --- Publication date ---
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# --- Publisher value ---
$meta.Range("B9").Value = "Alvearie Team"

# --- Replace the old "Contact" row content with the new Jurisdiction row ---
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Update root Extension element's Short/Definition text on the Elements sheet ---
$elem.Range("K2").Value = "Claim Adjustment Type"
$elem.Range("L2").Value = "The code for the claim's adjustment type"
